$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 18-33: rename "BD..." labels to "BDS..." and update their B/C/D values
# (the "BD" bulldozer/backtracking-style label group was renamed to "BDS" and
# its measured values were re-run, per "simplifikasi decition on v1 v2 v3").
$updates = @(
    @{ Row = 18; Name = "BDS";                 B = 160; C = 425; D = 1043 },
    @{ Row = 19; Name = "BDS-PPO";             B = 160; C = 425; D = 1043 },
    @{ Row = 20; Name = "BDS-TPF";             B = 113; C = 296; D = 653 },
    @{ Row = 21; Name = "BDS-PPO-TPF";         B = 113; C = 296; D = 653 },
    @{ Row = 22; Name = "BDS-BRC";             B = 103; C = 237; D = 588 },
    @{ Row = 23; Name = "BDS-PPO-BRC";         B = 103; C = 237; D = 588 },
    @{ Row = 24; Name = "BDS-BRC-TPF";         B = 118; C = 278; D = 733 },
    @{ Row = 25; Name = "BDS-PPO-BRC-TPF";     B = 118; C = 278; D = 733 },
    @{ Row = 26; Name = "BDS-GLM";             B = 111; C = 245; D = 697 },
    @{ Row = 27; Name = "BDS-PPO-GLM";         B = 111; C = 245; D = 697 },
    @{ Row = 28; Name = "BDS-GLM-TPF";         B = 113; C = 403; D = 803 },
    @{ Row = 29; Name = "BDS-PPO-GLM-TPF";     B = 113; C = 403; D = 803 },
    @{ Row = 30; Name = "BDS-BRC-GLM";         B = 164; C = 426; D = 1134 },
    @{ Row = 31; Name = "BDS-PPO-BRC-GLM";     B = 164; C = 426; D = 1134 },
    @{ Row = 32; Name = "BDS-BRC-GLM-TPF";     B = 109; C = 297; D = 724 },
    @{ Row = 33; Name = "BDS-PPO-BRC-GLM-TPF"; B = 109; C = 297; D = 724 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("A$r").Value = $u.Name
    $ws.Range("B$r").Value = $u.B
    $ws.Range("C$r").Value = $u.C
    $ws.Range("D$r").Value = $u.D
}

# Rows 50-65: only the label in column A changes from "BD..." to "BDS...";
# the B/C/D values are left untouched.
$nameOnly = @(
    @{ Row = 50; Name = "BDS-JPS" },
    @{ Row = 51; Name = "BDS-PPO-JPS" },
    @{ Row = 52; Name = "BDS-TPF-JPS" },
    @{ Row = 53; Name = "BDS-PPO-TPF-JPS" },
    @{ Row = 54; Name = "BDS-BRC-JPS" },
    @{ Row = 55; Name = "BDS-PPO-BRC-JPS" },
    @{ Row = 56; Name = "BDS-BRC-TPF-JPS" },
    @{ Row = 57; Name = "BDS-PPO-BRC-TPF-JPS" },
    @{ Row = 58; Name = "BDS-GLM-JPS" },
    @{ Row = 59; Name = "BDS-PPO-GLM-JPS" },
    @{ Row = 60; Name = "BDS-GLM-TPF-JPS" },
    @{ Row = 61; Name = "BDS-PPO-GLM-TPF-JPS" },
    @{ Row = 62; Name = "BDS-BRC-GLM-JPS" },
    @{ Row = 63; Name = "BDS-PPO-BRC-GLM-JPS" },
    @{ Row = 64; Name = "BDS-BRC-GLM-TPF-JPS" },
    @{ Row = 65; Name = "BDS-PPO-BRC-GLM-TPF-JPS" }
)

foreach ($u in $nameOnly) {
    $ws.Range("A$($u.Row)").Value = $u.Name
}

Write-Output "Updated $($updates.Count + $nameOnly.Count) rows"
